# Auto-generated edit script: apply 'Generate Report for Handoff' update
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Cells.Item(2, 1).Value = '69cdfbe9-67be-4967-9693-f158857d75ff.png'
$ws1.Cells.Item(2, 2).Value = 'Ready for handoff'
$ws1.Cells.Item(2, 3).Value = 'Ready for handoff'
$ws1.Cells.Item(2, 4).Value = '2016-03-22 00:56:59'
$ws1.Cells.Item(3, 1).Value = '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'
$ws1.Cells.Item(3, 2).Value = 'Ready for handoff'
$ws1.Cells.Item(3, 3).Value = 'Ready for handoff'
$ws1.Cells.Item(3, 4).Value = '2016-03-22 00:56:59'
$ws1.Cells.Item(4, 1).Value = 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png'
$ws1.Cells.Item(4, 2).Value = 'Ready for handoff'
$ws1.Cells.Item(4, 3).Value = 'Ready for handoff'
$ws1.Cells.Item(4, 4).Value = '2016-03-22 00:56:59'

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Cells.Item(2, 1).Value = '69cdfbe9-67be-4967-9693-f158857d75ff.png'
$ws2.Cells.Item(2, 2).Value = '.png'
$ws2.Cells.Item(2, 3).Value = 'Ready for handoff'
$ws2.Cells.Item(2, 4).Value = 'ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png'
$ws2.Cells.Item(2, 5).Value = '2016-03-22 00:56:55'
$ws2.Cells.Item(2, 8).Value = '0001-01-01 00:00:00'
$ws2.Cells.Item(2, 10).Value = 'IsDependency'
$ws2.Cells.Item(2, 11).Value = 'e2e\7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'
$ws2.Cells.Item(3, 1).Value = '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'
$ws2.Cells.Item(3, 2).Value = '.md'
$ws2.Cells.Item(3, 3).Value = 'Ready for handoff'
$ws2.Cells.Item(3, 4).Value = '7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.zh-cn.xlf'
$ws2.Cells.Item(3, 5).Value = '2016-03-22 00:56:55'
$ws2.Cells.Item(3, 8).Value = '0001-01-01 00:00:00'
$ws2.Cells.Item(3, 10).Value = 'Include'
$ws2.Cells.Item(4, 1).Value = 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png'
$ws2.Cells.Item(4, 2).Value = '.png'
$ws2.Cells.Item(4, 3).Value = 'Ready for handoff'
$ws2.Cells.Item(4, 4).Value = 'eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png'
$ws2.Cells.Item(4, 5).Value = '2016-03-22 00:56:55'
$ws2.Cells.Item(4, 8).Value = '0001-01-01 00:00:00'
$ws2.Cells.Item(4, 10).Value = 'IsDependency'
$ws2.Cells.Item(4, 11).Value = 'e2e\7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Cells.Item(2, 1).Value = '69cdfbe9-67be-4967-9693-f158857d75ff.png'
$ws3.Cells.Item(2, 2).Value = '.png'
$ws3.Cells.Item(2, 3).Value = 'Ready for handoff'
$ws3.Cells.Item(2, 4).Value = 'ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png'
$ws3.Cells.Item(2, 5).Value = '2016-03-22 00:56:59'
$ws3.Cells.Item(2, 8).Value = '0001-01-01 00:00:00'
$ws3.Cells.Item(2, 10).Value = 'IsDependency'
$ws3.Cells.Item(2, 11).Value = 'e2e\7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'
$ws3.Cells.Item(3, 1).Value = '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'
$ws3.Cells.Item(3, 2).Value = '.md'
$ws3.Cells.Item(3, 3).Value = 'Ready for handoff'
$ws3.Cells.Item(3, 4).Value = '7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.de-de.xlf'
$ws3.Cells.Item(3, 5).Value = '2016-03-22 00:56:59'
$ws3.Cells.Item(3, 8).Value = '0001-01-01 00:00:00'
$ws3.Cells.Item(3, 10).Value = 'Include'
$ws3.Cells.Item(4, 1).Value = 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png'
$ws3.Cells.Item(4, 2).Value = '.png'
$ws3.Cells.Item(4, 3).Value = 'Ready for handoff'
$ws3.Cells.Item(4, 4).Value = 'eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png'
$ws3.Cells.Item(4, 5).Value = '2016-03-22 00:56:59'
$ws3.Cells.Item(4, 8).Value = '0001-01-01 00:00:00'
$ws3.Cells.Item(4, 10).Value = 'IsDependency'
$ws3.Cells.Item(4, 11).Value = 'e2e\7bd6be59-01cb-44bd-a45e-b8b7c122553a.md'

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/69cdfbe9-67be-4967-9693-f158857d75ff.png', "", "", '69cdfbe9-67be-4967-9693-f158857d75ff.png')
$ws1.Hyperlinks.Add($ws1.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/7bd6be59-01cb-44bd-a45e-b8b7c122553a.md', "", "", '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md')
$ws1.Hyperlinks.Add($ws1.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/ed31b70d-5946-4f36-871a-51118a4ab7f6.png', "", "", 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png')

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/69cdfbe9-67be-4967-9693-f158857d75ff.png', "", "", '69cdfbe9-67be-4967-9693-f158857d75ff.png')
$ws2.Hyperlinks.Add($ws2.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e1ee8fe4055cd4539518604d4f19a2e85eb5d9c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png', "", "", 'ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png')
$ws2.Hyperlinks.Add($ws2.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/7bd6be59-01cb-44bd-a45e-b8b7c122553a.md', "", "", '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md')
$ws2.Hyperlinks.Add($ws2.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e1ee8fe4055cd4539518604d4f19a2e85eb5d9c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.zh-cn.xlf', "", "", '7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.zh-cn.xlf')
$ws2.Hyperlinks.Add($ws2.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/ed31b70d-5946-4f36-871a-51118a4ab7f6.png', "", "", 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png')
$ws2.Hyperlinks.Add($ws2.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e1ee8fe4055cd4539518604d4f19a2e85eb5d9c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png', "", "", 'eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png')

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/69cdfbe9-67be-4967-9693-f158857d75ff.png', "", "", '69cdfbe9-67be-4967-9693-f158857d75ff.png')
$ws3.Hyperlinks.Add($ws3.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a15ddd67ff521baadce86e81f596ea8f0f3c441a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png', "", "", 'ea62553b11f27c72ae2471c7eff19cf8e7c9815a.png')
$ws3.Hyperlinks.Add($ws3.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/7bd6be59-01cb-44bd-a45e-b8b7c122553a.md', "", "", '7bd6be59-01cb-44bd-a45e-b8b7c122553a.md')
$ws3.Hyperlinks.Add($ws3.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a15ddd67ff521baadce86e81f596ea8f0f3c441a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.de-de.xlf', "", "", '7bd6be59-01cb-44bd-a45e-b8b7c122553a.cb368a67e02e8c37df30e8f100333733143752b7.de-de.xlf')
$ws3.Hyperlinks.Add($ws3.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b3751bf45b20d5d749b3ce97b368d5b2ae873ae4/e2e/ed31b70d-5946-4f36-871a-51118a4ab7f6.png', "", "", 'ed31b70d-5946-4f36-871a-51118a4ab7f6.png')
$ws3.Hyperlinks.Add($ws3.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a15ddd67ff521baadce86e81f596ea8f0f3c441a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png', "", "", 'eeec4c3a394fed9960e3594c025d6b10f2e9ed57.png')

$ws1.Range('A2').Font.Underline = $true
$ws1.Range('A2').Font.Color = 15570276
$ws1.Range('A3').Font.Underline = $true
$ws1.Range('A3').Font.Color = 15570276
$ws1.Range('A4').Font.Underline = $true
$ws1.Range('A4').Font.Color = 15570276

$ws2.Range('A2').Font.Underline = $true
$ws2.Range('A2').Font.Color = 15570276
$ws2.Range('D2').Font.Underline = $true
$ws2.Range('D2').Font.Color = 15570276
$ws2.Range('A3').Font.Underline = $true
$ws2.Range('A3').Font.Color = 15570276
$ws2.Range('D3').Font.Underline = $true
$ws2.Range('D3').Font.Color = 15570276
$ws2.Range('A4').Font.Underline = $true
$ws2.Range('A4').Font.Color = 15570276
$ws2.Range('D4').Font.Underline = $true
$ws2.Range('D4').Font.Color = 15570276

$ws3.Range('A2').Font.Underline = $true
$ws3.Range('A2').Font.Color = 15570276
$ws3.Range('D2').Font.Underline = $true
$ws3.Range('D2').Font.Color = 15570276
$ws3.Range('A3').Font.Underline = $true
$ws3.Range('A3').Font.Color = 15570276
$ws3.Range('D3').Font.Underline = $true
$ws3.Range('D3').Font.Color = 15570276
$ws3.Range('A4').Font.Underline = $true
$ws3.Range('A4').Font.Color = 15570276
$ws3.Range('D4').Font.Underline = $true
$ws3.Range('D4').Font.Color = 15570276

$ws1.Range('D2').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Range('D3').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws1.Range('D4').NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range('E2').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range('H2').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range('E3').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range('H3').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range('E4').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range('H4').NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range('E2').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range('H2').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range('E3').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range('H3').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range('E4').NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range('H4').NumberFormat = "yyyy-mm-dd HH:mm:ss"

